$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 14.37161333333333
$ws.Range("N2").Value = 43.11484
$ws.Range("O2").Value = 0.4561705932627708
$ws.Range("P2").Value = 0.5019766122855294
$ws.Range("Q2").Value = 2.697269180937778
$ws.Range("R2").Value = 24.27542262844
$ws.Range("S2").Value = 0.4561705932627708
$ws.Range("T2").Value = 0.5019766122855294

# Row 3
$ws.Range("O3").Value = 0.01117178254830525
$ws.Range("P3").Value = 0.01229358849433434
$ws.Range("S3").Value = 0.01117178254830525
$ws.Range("T3").Value = 0.01229358849433434

# Row 4
$ws.Range("M4").Value = 2.798424666666667
$ws.Range("N4").Value = 8.395274000000001
$ws.Range("O4").Value = 0.08882503382091908
$ws.Range("P4").Value = 0.09774433122629669
$ws.Range("Q4").Value = 0.5252092742482223
$ws.Range("R4").Value = 4.726883468234001
$ws.Range("S4").Value = 0.08882503382091908
$ws.Range("T4").Value = 0.09774433122629669

# Row 5
$ws.Range("M5").Value = 8.624592
$ws.Range("N5").Value = 17.249184
$ws.Range("O5").Value = 0.2737539034788959
$ws.Range("P5").Value = 0.2008284606648142
$ws.Range("Q5").Value = 1.618666301424
$ws.Range("R5").Value = 9.711997808544
$ws.Range("S5").Value = 0.2737539034788959
$ws.Range("T5").Value = 0.2008284606648142

# Row 6
$ws.Range("M6").Value = 5.358313666666667
$ws.Range("N6").Value = 16.074941
$ws.Range("O6").Value = 0.1700786868891091
$ws.Range("P6").Value = 0.1871570073290255
$ws.Range("Q6").Value = 1.005650095064555
$ws.Range("R6").Value = 9.050850855581
$ws.Range("S6").Value = 0.1700786868891091
$ws.Range("T6").Value = 0.1871570073290255
